# Engimon.xlsx — "Spawn engimon liar, restyle starting menu, active engimon"
#
# 1. Column G ("Sprite Image") previously held a placeholder value
#    ("Characters/boy_stand_south.png") for every engimon row. Replace it
#    with the real per-engimon icon artwork (the first two rows get their
#    own distinct icon, the rest share the generic "icon011" icon).
# 2. Restyle the sheet view: clear the frozen/scrolled "topLeftCell" (was
#    A10) and move the active selection from F38 to G14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Icons/icon004.png"
$ws.Range("G3").Value = "Icons/icon333.png"
$ws.Range("G4:G37").Value = "Icons/icon011.png"

# Selecting G14 both updates the active selection and resets the
# sheet's scrolled top-left cell back to the sheet default (A1).
$ws.Range("G14").Select()
